$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 (Rule column) changes from "R40" to "1".
# The "1" must be stored as TEXT (not a number) while keeping the cell's
# existing style untouched. A plain `.Value = "1"` would auto-convert the
# numeric-looking string to a real number, so instead we write it as a text
# formula (a quoted string literal always evaluates to text in Excel) and
# then flatten the formula down to a literal value via copy/paste-values,
# which preserves the cell's type (text) without touching its style.
$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)
